# "Updates for appended spectroscopy"
# Replaces the (template) subject in row 2 of the Gas-Exchange summary sheet
# with the newly-appended spectroscopy subject Xe-037 (scanned/processed
# 2022-06-06 / 2022-06-17, Healthy-Cohort batch 20211119_HealthyCohort) and
# its associated numeric results, and re-applies the header/subject-row
# formatting plus a handful of narrower column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-apply the header row (row 1) & key/subject columns (row 2) text style ---
# (same "@" text format the sheet already uses for these label columns; Excel
# collapses it onto the identical existing style since there is no visible
# formatting change involved.)
$ws.Range("A1:BN1").NumberFormat = "@"
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"

# --- Row 2: swap the example subject for the newly appended one ---
$ws.Range("A2").Value = "Xe-037"
$ws.Range("B2").Value = "2022-06-06"
$ws.Range("E2").Value = "2022-06-17"
$ws.Range("F2").Value = "20211119_HealthyCohort"
# C2 (TE90) and D2 (Flip_Angle) are unchanged.

# --- Row 2: updated spectroscopy results ---
$ws.Range("G2").Value = -7396.0699979126648
$ws.Range("H2").Value = -686.60554883478142
$ws.Range("I2").Value = -0.21277834592285871
$ws.Range("J2").Value = 12.452751020758702
$ws.Range("K2").Value = 1.7535812844099699
$ws.Range("L2").Value = 0.99421228726053779
$ws.Range("M2").Value = 1.2069861407667319
$ws.Range("N2").Value = 22.869523944135615
$ws.Range("O2").Value = -162.64375004412031
$ws.Range("P2").Value = -73.641403860233098
$ws.Range("Q2").Value = 0.50422057628367611
$ws.Range("R2").Value = 1.0824047123364617
$ws.Range("S2").Value = 116.16012608389605
$ws.Range("T2").Value = 18.763269575972821
$ws.Range("U2").Value = 42.522362885128132
$ws.Range("V2").Value = 40.173244040160888
$ws.Range("W2").Value = 24.592538019237868
$ws.Range("X2").Value = 12.565206315223557
$ws.Range("Y2").Value = 0.58971121836490292
$ws.Range("Z2").Value = 0.21720213227423599
$ws.Range("AA2").Value = 0.96565968210125641
$ws.Range("AB2").Value = 0.3059578572291124
$ws.Range("AC2").Value = 0.83424191431792716
$ws.Range("AD2").Value = 0.28127478092069452
$ws.Range("AE2").Value = 0.46028858435511455
$ws.Range("AF2").Value = 0.17206644731132772
$ws.Range("AG2").Value = 0.56578582033089164
$ws.Range("AH2").Value = 0.18373352569072537
$ws.Range("AI2").Value = 2.8564916812276837
$ws.Range("AJ2").Value = 5.5230570538631012
$ws.Range("AK2").Value = 17.997417002203147
$ws.Range("AL2").Value = 38.965281470789328
$ws.Range("AM2").Value = 31.056749981007371
$ws.Range("AN2").Value = 3.601002810909367
$ws.Range("AO2").Value = 0.6881989520606866
$ws.Range("AP2").Value = 10.581058887933056
$ws.Range("AQ2").Value = 33.189958551654023
$ws.Range("AR2").Value = 34.089309454915153
$ws.Range("AS2").Value = 18.033940721044811
$ws.Range("AT2").Value = 3.4175334323922733
$ws.Range("AU2").Value = 0.20333150856338467
$ws.Range("AV2").Value = 8.0081332603425359
$ws.Range("AW2").Value = 38.539141315398453
$ws.Range("AX2").Value = 33.205599436928132
$ws.Range("AY2").Value = 14.983968092594042
$ws.Range("AZ2").Value = 4.3090638930163445
$ws.Range("BA2").Value = 0.54743098459372808
$ws.Range("BB2").Value = 0.20333150856338467
$ws.Range("BC2").Value = 0.84460780480175179
$ws.Range("BD2").Value = 7.6014702432157657
$ws.Range("BE2").Value = 29.00602174083053
$ws.Range("BF2").Value = 50.066473762414951
$ws.Range("BG2").Value = 11.855791037772738
$ws.Range("BH2").Value = 0.62563541096426056
$ws.Range("BI2").Value = 0.71948072260889961
$ws.Range("BJ2").Value = 5.0207241729881913
$ws.Range("BK2").Value = 34.034566356455777
$ws.Range("BL2").Value = 56.541800265895048
$ws.Range("BM2").Value = 3.6443262688668181
$ws.Range("BN2").Value = 0.039102213185266285

# --- Narrow a handful of columns to fit the appended data better ---
# (ColumnWidth is quantized to whole pixels by Excel; these are the
# closest attainable widths to the target layout.)
$ws.Columns.Item(2).ColumnWidth = 9.5                  # B  -> ~10.29
$ws.Columns.Item(6).ColumnWidth = 14.166666666666666   # F  -> 15
$ws.Columns.Item(10).ColumnWidth = 9.666666666666666   # J  -> ~10.57
$ws.Columns.Item(13).ColumnWidth = 10.0                # M  -> ~10.86
$ws.Columns.Item(14).ColumnWidth = 9.666666666666666   # N  -> ~10.57
$ws.Columns.Item(16).ColumnWidth = 10.0                # P  -> ~10.86
$ws.Columns.Item(21).ColumnWidth = 8.0                 # U  -> ~8.86
$ws.Columns.Item(24).ColumnWidth = 8.333333333333334   # X  -> ~9.14
